# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose automatic dialog-act annotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 13;  DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 24;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 46;  DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 51;  DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 63;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 67;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 79;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 80;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 88;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 89;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 93;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 98;  DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 100; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 102; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 104; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 109; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 110; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 111; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 131; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 134; DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 144; DAMSLTag = "%";  DialogAct = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
